$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Basketball Game Sales")

# New "Calories" column (F) on the Basketball Game Sales sheet.
$ws.Range("F1").Value = "Calories"

for ($r = 2; $r -le 200; $r++) {
    $ws.Range("F$r").Formula = "=VLOOKUP(A$r,Calories!`$A`$1:`$B`$15,2,FALSE)"
}

# Make "Basketball Game Sales" the active sheet/tab, with the selection on I6,
# leaving the Calories sheet's own selection (D40) untouched.
$ws.Select() | Out-Null
$ws.Range("I6").Select() | Out-Null
